# pres poll 20201011 clean data and data dictionary update
# Add data-dictionary rows describing the new pres_poll columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pres_poll")

# Column A ("Feature" names) filled top to bottom first.
$ws.Range("A2").Value = "Day"
$ws.Range("A3").Value = "Len"
$ws.Range("A4").Value = "State"
$ws.Range("A5").Value = "EV"
$ws.Range("A6").Value = "Dem"
$ws.Range("A7").Value = "GOP"
$ws.Range("A8").Value = "Ind"
$ws.Range("A9").Value = "Date"
$ws.Range("A10").Value = "Pollster"

# Column B ("Description" text) filled top to bottom next.
$ws.Range("B2").Value = "Day of the year at the midpoint of the poll - so a poll from Feb 02 to Feb 04 would have Feb 03 as the midpoint and this is day 34"
$ws.Range("B3").Value = "Duration of Poll"
$ws.Range("B4").Value = "State"
$ws.Range("B5").Value = "Electoral Vote"
$ws.Range("B6").Value = "Democrat Party %"
$ws.Range("B7").Value = "Republican Party %"
$ws.Range("B8").Value = "Independent Party %"
$ws.Range("B9").Value = "End-date of data collection for poll"
$ws.Range("B10").Value = "Pollster"

# Leave pres_poll as the active/selected sheet with B10 the active cell,
# matching the saved view state.
$ws.Activate()
[void]$ws.Range("B10").Select()
